$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 15 "data" rows (rows 1-15) followed by a final
# trailer row (row 16). We need to append two more copies of the 15 data
# rows before the trailer row, so the trailer ends up at row 46.

$blockSize = 15
$numCols = 5

# Capture the block of 15 rows (A1:E15) as values, cell by cell.
$blockValues = @{}
for ($r = 1; $r -le $blockSize; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $blockValues["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Capture the trailer row (row 16) values before we move anything.
$lastRowIndex = 16
$trailerValues = @{}
for ($c = 1; $c -le $numCols; $c++) {
    $trailerValues[$c] = $ws.Cells.Item($lastRowIndex, $c).Value2
}

# Write two more copies of the block starting right after the existing data
# (append-to-writer style).
for ($copy = 1; $copy -le 2; $copy++) {
    $destStart = $blockSize * $copy + 1
    for ($r = 1; $r -le $blockSize; $r++) {
        $destRow = $destStart + $r - 1
        for ($c = 1; $c -le $numCols; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $blockValues["$r,$c"]
        }
    }
}

# Finally, write the trailer row at its new location (after all copies).
$newTrailerRow = $blockSize * 3 + 1
for ($c = 1; $c -le $numCols; $c++) {
    $ws.Cells.Item($newTrailerRow, $c).Value = $trailerValues[$c]
}
